$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Expand the "History" table from 1 column to 6 columns (A1:F2) ---
$lo = $ws.ListObjects.Item("History")
$lo.Resize($ws.Range("A1:F2"))

# Renaming columns through the header-row cell values (drives both the
# sharedStrings table and the table's tableColumn/@name).
$ws.Range("A1").Value = "Timestamp"
$ws.Range("B1").Value = "Context"
$ws.Range("C1").Value = "Title"
$ws.Range("D1").Value = "Host"
$ws.Range("E1").Value = "URL"
$ws.Range("F1").Value = "User Agent"

# --- Column widths for the new columns ---
$ws.Columns.Item(1).ColumnWidth = 19.833333333333332
$ws.Columns.Item(2).ColumnWidth = 18
$ws.Columns.Item(3).ColumnWidth = 32.166666666666664
$ws.Columns.Item(4).ColumnWidth = 15.333333333333334
$ws.Columns.Item(5).ColumnWidth = 12.833333333333334
$ws.Columns.Item(6).ColumnWidth = 16.666666666666668

# --- Selection moves to B4 ---
$ws.Range("B4").Select() | Out-Null
